{"js": "// Update the worksheet's 25 \"two-digit \u00f7 one-digit\" answer cells.\n// The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16)\n// carries data, the rows in between are blank spacer rows. Each old answer\n// string is unique in the document, so we locate cells by their current\n// (pre-edit) text and overwrite them with the new answer text in place,\n// preserving the run/paragraph formatting already on the cell.\nconst replacements = [\n  [\"55\u00f74=13, 3\", \"43\u00f77=6, 1\"],\n  [\"97\u00f75=19, 2\", \"70\u00f75=14, 0\"],\n  [\"82\u00f74=20, 2\", \"64\u00f77=9, 1\"],\n  [\"89\u00f73=29, 2\", \"10\u00f73=3, 1\"],\n  [\"62\u00f74=15, 2\", \"59\u00f79=6, 5\"],\n  [\"84\u00f76=14, 0\", \"30\u00f74=7, 2\"],\n  [\"22\u00f79=2, 4\", \"34\u00f75=6, 4\"],\n  [\"31\u00f79=3, 4\", \"43\u00f77=6, 1\"],\n  [\"11\u00f77=1, 4\", \"67\u00f75=13, 2\"],\n  [\"82\u00f72=41, 0\", \"75\u00f72=37, 1\"],\n  [\"54\u00f78=6, 6\", \"47\u00f74=11, 3\"],\n  [\"22\u00f73=7, 1\", \"15\u00f72=7, 1\"],\n  [\"86\u00f78=10, 6\", \"81\u00f75=16, 1\"],\n  [\"46\u00f72=23, 0\", \"13\u00f73=4, 1\"],\n  [\"57\u00f73=19, 0\", \"26\u00f76=4, 2\"],\n  [\"40\u00f77=5, 5\", \"16\u00f72=8, 0\"],\n  [\"71\u00f74=17, 3\", \"29\u00f74=7, 1\"],\n  [\"52\u00f75=10, 2\", \"81\u00f74=20, 1\"],\n  [\"13\u00f79=1, 4\", \"95\u00f77=13, 4\"],\n  [\"23\u00f77=3, 2\", \"71\u00f75=14, 1\"],\n  [\"98\u00f78=12, 2\", \"41\u00f77=5, 6\"],\n  [\"34\u00f79=3, 7\", \"48\u00f72=24, 0\"],\n  [\"71\u00f72=35, 1\", \"37\u00f76=6, 1\"],\n  [\"43\u00f73=14, 1\", \"29\u00f78=3, 5\"],\n  [\"74\u00f78=9, 2\", \"62\u00f77=8, 6\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet's 25 \"two-digit \u00f7 one-digit\" answer cells.\n# The table has 20 rows x 5 columns; only every 4th row (0, 4, 8, 12, 16)\n# carries data, the rows in between are blank spacer rows. Each old answer\n# string is unique in the document, so Find/Replace on the whole document\n# content locates each cell's run and rewrites its text in place, preserving\n# the run/paragraph formatting already on the cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"55\u00f74=13, 3\", \"43\u00f77=6, 1\"),\n  @(\"97\u00f75=19, 2\", \"70\u00f75=14, 0\"),\n  @(\"82\u00f74=20, 2\", \"64\u00f77=9, 1\"),\n  @(\"89\u00f73=29, 2\", \"10\u00f73=3, 1\"),\n  @(\"62\u00f74=15, 2\", \"59\u00f79=6, 5\"),\n  @(\"84\u00f76=14, 0\", \"30\u00f74=7, 2\"),\n  @(\"22\u00f79=2, 4\", \"34\u00f75=6, 4\"),\n  @(\"31\u00f79=3, 4\", \"43\u00f77=6, 1\"),\n  @(\"11\u00f77=1, 4\", \"67\u00f75=13, 2\"),\n  @(\"82\u00f72=41, 0\", \"75\u00f72=37, 1\"),\n  @(\"54\u00f78=6, 6\", \"47\u00f74=11, 3\"),\n  @(\"22\u00f73=7, 1\", \"15\u00f72=7, 1\"),\n  @(\"86\u00f78=10, 6\", \"81\u00f75=16, 1\"),\n  @(\"46\u00f72=23, 0\", \"13\u00f73=4, 1\"),\n  @(\"57\u00f73=19, 0\", \"26\u00f76=4, 2\"),\n  @(\"40\u00f77=5, 5\", \"16\u00f72=8, 0\"),\n  @(\"71\u00f74=17, 3\", \"29\u00f74=7, 1\"),\n  @(\"52\u00f75=10, 2\", \"81\u00f74=20, 1\"),\n  @(\"13\u00f79=1, 4\", \"95\u00f77=13, 4\"),\n  @(\"23\u00f77=3, 2\", \"71\u00f75=14, 1\"),\n  @(\"98\u00f78=12, 2\", \"41\u00f77=5, 6\"),\n  @(\"34\u00f79=3, 7\", \"48\u00f72=24, 0\"),\n  @(\"71\u00f72=35, 1\", \"37\u00f76=6, 1\"),\n  @(\"43\u00f73=14, 1\", \"29\u00f78=3, 5\"),\n  @(\"74\u00f78=9, 2\", \"62\u00f77=8, 6\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n"}
